$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the private key values from A2 and A3 (sensitive data removed from the sheet)
$ws.Range("A2:A3").ClearContents()

# Move the active selection to A2, matching the saved state
$ws.Range("A2").Select()
